# Update the 20x5 table of arithmetic problems ("within100.docx")
# to the new set of expressions, cell by cell, preserving formatting.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "4+8="
$t.Cell(1,2).Range.Text = "81-39="
$t.Cell(1,3).Range.Text = "62-60="
$t.Cell(1,4).Range.Text = "91-5="
$t.Cell(1,5).Range.Text = "12+31="
$t.Cell(2,1).Range.Text = "2+77="
$t.Cell(2,2).Range.Text = "68-30="
$t.Cell(2,3).Range.Text = "52+43="
$t.Cell(2,4).Range.Text = "75-13="
$t.Cell(2,5).Range.Text = "33+62="
$t.Cell(3,1).Range.Text = "50+11="
$t.Cell(3,2).Range.Text = "38+37="
$t.Cell(3,3).Range.Text = "37-36="
$t.Cell(3,4).Range.Text = "60-27="
$t.Cell(3,5).Range.Text = "2+5="
$t.Cell(4,1).Range.Text = "93+3="
$t.Cell(4,2).Range.Text = "67+8="
$t.Cell(4,3).Range.Text = "34+55="
$t.Cell(4,4).Range.Text = "83-78="
$t.Cell(4,5).Range.Text = "14+15="
$t.Cell(5,1).Range.Text = "64-48="
$t.Cell(5,2).Range.Text = "93-54="
$t.Cell(5,3).Range.Text = "61-14="
$t.Cell(5,4).Range.Text = "80-36="
$t.Cell(5,5).Range.Text = "78+9="
$t.Cell(6,1).Range.Text = "69-69="
$t.Cell(6,2).Range.Text = "48+7="
$t.Cell(6,3).Range.Text = "0+93="
$t.Cell(6,4).Range.Text = "15+45="
$t.Cell(6,5).Range.Text = "15+57="
$t.Cell(7,1).Range.Text = "92-9="
$t.Cell(7,2).Range.Text = "36+50="
$t.Cell(7,3).Range.Text = "18+54="
$t.Cell(7,4).Range.Text = "58+37="
$t.Cell(7,5).Range.Text = "25+35="
$t.Cell(8,1).Range.Text = "36+32="
$t.Cell(8,2).Range.Text = "91-54="
$t.Cell(8,3).Range.Text = "18+57="
$t.Cell(8,4).Range.Text = "51+10="
$t.Cell(8,5).Range.Text = "80-13="
$t.Cell(9,1).Range.Text = "75-17="
$t.Cell(9,2).Range.Text = "84-15="
$t.Cell(9,3).Range.Text = "6+83="
$t.Cell(9,4).Range.Text = "5-1="
$t.Cell(9,5).Range.Text = "86-18="
$t.Cell(10,1).Range.Text = "17+57="
$t.Cell(10,2).Range.Text = "62-40="
$t.Cell(10,3).Range.Text = "30+47="
$t.Cell(10,4).Range.Text = "71+1="
$t.Cell(10,5).Range.Text = "49+15="
$t.Cell(11,1).Range.Text = "64+10="
$t.Cell(11,2).Range.Text = "72-71="
$t.Cell(11,3).Range.Text = "46+2="
$t.Cell(11,4).Range.Text = "7+9="
$t.Cell(11,5).Range.Text = "17+78="
$t.Cell(12,1).Range.Text = "18+44="
$t.Cell(12,2).Range.Text = "40-18="
$t.Cell(12,3).Range.Text = "97-37="
$t.Cell(12,4).Range.Text = "56+38="
$t.Cell(12,5).Range.Text = "3+80="
$t.Cell(13,1).Range.Text = "18+22="
$t.Cell(13,2).Range.Text = "8+73="
$t.Cell(13,3).Range.Text = "3+61="
$t.Cell(13,4).Range.Text = "41+24="
$t.Cell(13,5).Range.Text = "91+2="
$t.Cell(14,1).Range.Text = "37+25="
$t.Cell(14,2).Range.Text = "88-14="
$t.Cell(14,3).Range.Text = "38+46="
$t.Cell(14,4).Range.Text = "90-14="
$t.Cell(14,5).Range.Text = "82-55="
$t.Cell(15,1).Range.Text = "82+0="
$t.Cell(15,2).Range.Text = "93-73="
$t.Cell(15,3).Range.Text = "95-58="
$t.Cell(15,4).Range.Text = "75+0="
$t.Cell(15,5).Range.Text = "27-20="
$t.Cell(16,1).Range.Text = "26+69="
$t.Cell(16,2).Range.Text = "55-10="
$t.Cell(16,3).Range.Text = "48-37="
$t.Cell(16,4).Range.Text = "96-0="
$t.Cell(16,5).Range.Text = "56+29="
$t.Cell(17,1).Range.Text = "56+11="
$t.Cell(17,2).Range.Text = "1+11="
$t.Cell(17,3).Range.Text = "84-47="
$t.Cell(17,4).Range.Text = "31-13="
$t.Cell(17,5).Range.Text = "57+21="
$t.Cell(18,1).Range.Text = "45-12="
$t.Cell(18,2).Range.Text = "36+47="
$t.Cell(18,3).Range.Text = "68-0="
$t.Cell(18,4).Range.Text = "86-46="
$t.Cell(18,5).Range.Text = "89-32="
$t.Cell(19,1).Range.Text = "7+44="
$t.Cell(19,2).Range.Text = "52+22="
$t.Cell(19,3).Range.Text = "76-13="
$t.Cell(19,4).Range.Text = "94-92="
$t.Cell(19,5).Range.Text = "99-69="
$t.Cell(20,1).Range.Text = "15-5="
$t.Cell(20,2).Range.Text = "8+87="
$t.Cell(20,3).Range.Text = "39+23="
$t.Cell(20,4).Range.Text = "79-17="
$t.Cell(20,5).Range.Text = "16-4="
